$p = $ppt.ActivePresentation

# Commit: "inversão de slides: scratch e hour of code"
# Slide 4 (Scratch) and Slide 5 (Hora do Código) swap places.
$s = $p.Slides.Item(4)
$s.MoveTo(5)
